# ValidMultiplePageData.xlsx - "Whole Lot Of Changes"
#
# 1) Rename the fully-qualified test-class strings stored in the shared
#    string table (ExcelFileMultiplePageTest -> ExcelFileProcessorMultiplePageTest)
#    on the "Data" sheet. Using Cells.Replace (instead of rewriting each
#    cell's .Value) mutates the shared-string entries in place so every
#    cell that referenced them keeps pointing at the same shared-string
#    index, exactly like the authored diff.
# 2) Move the selection on the "Data" sheet from A9 to A8.
# 3) Best-effort: nudge the workbook window's on-screen position/size to
#    match the recorded workbookView bounds.

$wb = $excel.ActiveWorkbook

$dataSheet = $excel.Worksheets.Item("Data")

# --- 1) Shared-string rename (FirstMultiplePage / ThirdMultiplePage / SecondMultiplePage rows) ---
$dataSheet.Cells.Replace("ExcelFileMultiplePageTest", "ExcelFileProcessorMultiplePageTest")

# --- 2) Selection moves from A9 to A8 on the already-active "Data" tab ---
$dataSheet.Range("A8").Select()

# --- 3) Best-effort window geometry (xWindow/yWindow/windowWidth) ---
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 2190
$win.Width = 14490
$win.Height = 6750
